$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed as a number by Excel's
# type inference; format as text first so the literal string is preserved,
# matching the source inlineStr cell content exactly.
$textCells = @(
    'D6','D7','D8','D9','D10','D12','D13','D14','D16','D19','D20','D21','D22','D24','D26','D27','D28','D29','D30','D31','D33','D35','D36','D37','D38','D39','D40','D41','D43','D44','D45','D46','D48','D49','D50'
)
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$values = @{
    'D2' = '23.097.30'
    'D3' = '1.600.45'
    'D6' = '300.95'
    'D7' = '0.3776'
    'D8' = '0.3634'
    'D9' = '50.02'
    'D10' = '1.254'
    'D12' = '0.08123'
    'D13' = '22.98'
    'D14' = '6.562'
    'D16' = '7.350'
    'D17' = '1.599.48'
    'D19' = '0.06864'
    'D20' = '18.26'
    'D21' = '6.538'
    'D22' = '0.5565'
    'D24' = '12.95'
    'D25' = '23.089.81'
    'D26' = '2.337'
    'D27' = '2.712'
    'D28' = '21.04'
    'D29' = '149.59'
    'D30' = '5.265'
    'D31' = '131.73'
    'D33' = '6.827'
    'D34' = '1.775.92'
    'D35' = '0.9513'
    'D36' = '0.07641'
    'D37' = '0.02723'
    'D38' = '0.2542'
    'D39' = '6.188'
    'D40' = '0.08883'
    'D41' = '10.01'
    'D43' = '0.7059'
    'D44' = '12.60'
    'D45' = '15.19'
    'D46' = '0.6565'
    'D48' = '2.294'
    'D49' = '3.975'
    'D50' = '131.97'
    'E2' = '  -3.61%  '
    'E3' = '  -2.91%  '
    'E4' = '  +0.05%  '
    'E5' = '  +0.03%  '
    'E6' = '  -2.98%  '
    'E7' = '  -2.82%  '
    'E8' = '  -4.48%  '
    'E9' = '  -3.77%  '
    'E10' = '  -6.73%  '
    'E11' = '  -0.03%  '
    'E12' = '  -3.68%  '
    'E13' = '  -3.80%  '
    'E14' = '  -6.93%  '
    'E15' = '  -4.07%  '
    'E16' = '  -8.65%  '
    'E17' = '  -3.17%  '
    'E18' = '  -2.69%  '
    'E19' = '  -1.78%  '
    'E20' = '  -7.01%  '
    'E21' = '  -5.83%  '
    'E22' = '  -6.53%  '
    'E23' = '  +0.05%  '
    'E24' = '  -5.74%  '
    'E25' = '  -3.59%  '
    'E26' = '  -5.05%  '
    'E27' = '  -8.02%  '
    'E28' = '  -4.49%  '
    'E29' = '  -2.49%  '
    'E30' = '  -2.35%  '
    'E31' = '  -4.61%  '
    'E32' = '  -3.43%  '
    'E33' = '  -12.60%  '
    'E34' = '  -2.95%  '
    'E35' = '  -6.08%  '
    'E36' = '  -5.74%  '
    'E37' = '  -7.04%  '
    'E38' = '  -4.96%  '
    'E39' = '  -8.27%  '
    'E40' = '  -2.11%  '
    'E41' = '  -6.43%  '
    'E42' = '  -3.08%  '
    'E43' = '  -6.85%  '
    'E44' = '  -5.86%  '
    'E45' = '  -6.21%  '
    'E46' = '  -5.43%  '
    'E48' = '  -6.27%  '
    'E49' = '  -2.87%  '
    'E50' = '  -1.56%  '
    'E51' = '  -4.40%  '
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
